$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (interested count) for two events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3206
$ws1.Range("F6").Value = 137

# Sheet "全部类型" - same two events duplicated here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3206
$ws4.Range("F11").Value = 137
